# Adds two new verb rows' missing conjugations (nehmen / schlafen), a new
# "category" column to the Table1 listobject, and a new row for the verb
# "geben" (to give), per the "verbs added to resources" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$lo = $ws.ListObjects.Item(1)

# ---------------------------------------------------------------------
# 1. Fill in the missing conjugation/IPA cells for two existing verbs.
#    (These cells already exist in the sheet with style "2" applied -
#    D28/E28 already held values; we only add the ones that were blank.)
# ---------------------------------------------------------------------

# Row 28: nehmen / to take -- add IPA + er/sie/es, wir, ihr, Sie, sie
$ws.Cells.Item(28, 3).Value  = "ˈneːmən"   # C28 IPA
$ws.Cells.Item(28, 6).Value  = "nimmt"     # F28 er/sie/es
$ws.Cells.Item(28, 7).Value  = "nehmen"    # G28 wir
$ws.Cells.Item(28, 8).Value  = "nehmt"     # H28 ihr
$ws.Cells.Item(28, 9).Value  = "nehmen"    # I28 Sie
$ws.Cells.Item(28, 10).Value = "nehmen"    # J28 sie

# Row 34: schlafen / sleep -- add IPA + every conjugation
$ws.Cells.Item(34, 3).Value  = "ˈʃlaːfn"   # C34 IPA
$ws.Cells.Item(34, 4).Value  = "schlafe"   # D34 ich
$ws.Cells.Item(34, 5).Value  = "schläfst"  # E34 du
$ws.Cells.Item(34, 6).Value  = "schläft"   # F34 er/sie/es
$ws.Cells.Item(34, 7).Value  = "schlafen"  # G34 wir
$ws.Cells.Item(34, 8).Value  = "schlaft"   # H34 ihr
$ws.Cells.Item(34, 9).Value  = "schlafen"  # I34 Sie
$ws.Cells.Item(34, 10).Value = "schlafen"  # J34 sie

# ---------------------------------------------------------------------
# 2. Add a new "category" table column (becomes column K) and a new
#    table row (becomes row 52) for the verb "geben".
# ---------------------------------------------------------------------

$newCol = $lo.ListColumns().Add()
$ws.Cells.Item(1, 11).Value = "category"

$newRow = $lo.ListRows().Add()

# ---------------------------------------------------------------------
# 3. Stamp the plain data style (as used throughout columns D:J, e.g.
#    D4) onto the newly-created column/row cells so the blank ones
#    match the rest of the table's formatting.
# ---------------------------------------------------------------------

$styleDonor = $ws.Range("D4")
$styleDonor.Copy()
$ws.Range("K1:K51").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("D52:K52").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Re-apply the header text (PasteSpecial only touched formatting, but
# make sure it's still correct) and fill in the new row's data.
$ws.Cells.Item(1, 11).Value = "category"

$ws.Cells.Item(52, 1).Value  = "geben"       # verb
$ws.Cells.Item(52, 2).Value  = "to give"     # meaning
$ws.Cells.Item(52, 3).Value  = "ˈɡeːbn"      # IPA
$ws.Cells.Item(52, 4).Value  = "gebe"        # ich
$ws.Cells.Item(52, 5).Value  = "gibst"       # du
$ws.Cells.Item(52, 6).Value  = "gibt"        # er/sie/es
$ws.Cells.Item(52, 7).Value  = "geben"       # wir
$ws.Cells.Item(52, 8).Value  = "gebt"        # ihr
$ws.Cells.Item(52, 9).Value  = "geben"       # Sie
$ws.Cells.Item(52, 10).Value = "geben"       # sie
$ws.Cells.Item(52, 11).Value = "transitive"  # category

# ---------------------------------------------------------------------
# 4. Mirror the selection change recorded in the commit (cosmetic).
# ---------------------------------------------------------------------

$ws.Range("F28").Select()
